$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B35").Copy()
$ws.Range("B38").PasteSpecial(-4122)
$ws.Range("C35").Copy()
$ws.Range("C38").PasteSpecial(-4122)
$ws.Range("B38").Value = "Mise en place propre modele MVC"
$ws.Range("D38").Value = 4
